$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 2 through 27, keeping only the header row (row 1)
$ws.Range("A2:B27").EntireRow.Delete()
